# Split the single "8.1.2024." run in the "Đurđevac, 8.1.2024." paragraph
# into six runs:
#   {{ dan_naloga }} | . |  {{ mj_naloga }} | . |   | 2024.
# All six runs keep the exact original run formatting
# (rStyle "s9", Arial/Arial/Arial fonts, color 000000, sz 14, szCs 14).

$d = $word.ActiveDocument

# Locate the run that currently holds the literal date "8.1.2024."
$target = $d.Content
$target.Find.Execute("8.1.2024.", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "", 0)
$start = $target.Start

# A nearby run ("URBROJ:2137-37-24-1") carries the identical character
# formatting (same rStyle/fonts/color/size) and is left untouched by this
# edit, so it is used as a formatting template for the new runs.
$tmplRng = $d.Content
$tmplRng.Find.Execute("URBROJ:2137-37-24-1", $true, $false, $false, $false, `
                       $false, $true, 1, $false, "", 0)
$tmplLen = $tmplRng.End - $tmplRng.Start

# Remove the old "8.1.2024." text; $start now marks the insertion point.
$target.Text = ""

# The six replacement pieces, concatenating back to:
#   {{ dan_naloga }}. {{ mj_naloga }}. 2024.
$pieces = @(
    "{{ dan_naloga }}",
    ".",
    " {{ mj_naloga }}",
    ".",
    " ",
    "2024."
)

$pos = $start
foreach ($piece in $pieces) {
    # Grab a fresh copy of the template run's formatting for every piece
    # (the reference must be re-acquired each time, since earlier copies
    # become stale once the document text shifts).
    $tmpl = $d.Content
    $tmpl.Find.Execute("URBROJ:2137-37-24-1", $true, $false, $false, $false, `
                        $false, $true, 1, $false, "", 0)
    $tmpl.Copy()

    $dest = $d.Range($pos, $pos)
    $dest.Paste()

    $newRun = $d.Range($pos, $pos + $tmplLen)
    # Nudge a formatting property away from the template's value first so
    # that the upcoming Text assignment is not silently coalesced into the
    # identically-formatted run sitting right before it.
    $newRun.Font.Color = 255
    $newRun.Text = $piece

    # Restore the correct (black) color now that the run holds its final
    # text; this property-only change does not trigger run coalescing.
    $fixRun = $d.Range($pos, $pos + $piece.Length)
    $fixRun.Font.Color = 0

    $pos = $pos + $piece.Length
}
